$d = $word.ActiveDocument

$replacements = @(
    @("2025-10-01 Wednesday", "2025-10-02 Thursday"),
    @("23×64=1472", "69×40=2760"),
    @("83×53=4399", "97×60=5820"),
    @("69×28=1932", "72×29=2088"),
    @("26×21=546", "77×91=7007"),
    @("82×42=3444", "71×40=2840"),
    @("78×46=3588", "45×60=2700"),
    @("67×47=3149", "18×21=378"),
    @("94×59=5546", "40×95=3800"),
    @("59×31=1829", "33×93=3069"),
    @("21×71=1491", "22×39=858"),
    @("52×65=3380", "55×39=2145"),
    @("74×55=4070", "43×17=731"),
    @("32×42=1344", "67×69=4623"),
    @("58×86=4988", "94×85=7990"),
    @("49×16=784", "98×54=5292"),
    @("76×16=1216", "29×61=1769"),
    @("61×61=3721", "40×75=3000"),
    @("48×96=4608", "39×63=2457"),
    @("40×91=3640", "48×71=3408"),
    @("21×80=1680", "17×63=1071"),
    @("80×47=3760", "93×24=2232"),
    @("35×77=2695", "73×87=6351"),
    @("78×51=3978", "60×45=2700"),
    @("77×65=5005", "29×55=1595"),
    @("99×47=4653", "34×76=2584")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
